$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data (row 6) appended below the existing table (rows 1-5)
$ws.Cells.Item(6, 1).Value = 42607.886504629627
$ws.Cells.Item(6, 2).Value = 2
$ws.Cells.Item(6, 3).Value = 58
$ws.Cells.Item(6, 4).Value = 39
$ws.Cells.Item(6, 5).Value = 36
$ws.Cells.Item(6, 6).Value = 63
$ws.Cells.Item(6, 7).Value = 20968
$ws.Cells.Item(6, 8).Value = 16680
$ws.Cells.Item(6, 9).Value = 2798
$ws.Cells.Item(6, 10).Value = 363
$ws.Cells.Item(6, 11).Value = 244
$ws.Cells.Item(6, 12).Value = 11
$ws.Cells.Item(6, 13).Value = 19
$ws.Cells.Item(6, 14).Value = "Noun"

# Match the date/time number format used by the other cells in column A
$ws.Cells.Item(6, 1).NumberFormat = "m/d/yy h:mm"
